# Rogue One TrainSystem schedule.xlsx edit
# Commit: "GUI now displays Red Line schedule generated from Excel file"
#
# Summary of change:
#  - Rename sheet "Train" -> "Red"
#  - Add a new sheet "Green" right after "Red"
#  - Rebuild both "Red" and "Green" sheets with a driver/departure/station
#    schedule table (replacing the old train-telemetry column headers)
#  - "Green" additionally carries a time-of-day schedule computed with
#    TIME() offset formulas
#  - Personnel sheet is untouched content-wise

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Sheets: rename + add
# ---------------------------------------------------------------------
$red = $wb.Worksheets.Item("Train")
$red.Name = "Red"

$green = $wb.Worksheets.Add($null, $red)
$green.Name = "Green"

# ---------------------------------------------------------------------
# 2. "Red" sheet — header row
# ---------------------------------------------------------------------
$red.Range("A1").Value = "Train ID"
$red.Range("B1").Value = "DRIVER"
$red.Range("C1").Value = "DEPARTURE"
$red.Range("D1").Value = "SHADYSIDE"
$red.Range("E1").Value = "HERRON"
$red.Range("F1").Value = "SWISSVILLE"
$red.Range("G1").Value = "PENN STATION"
$red.Range("H1").Value = "STEEL PLAZA"
$red.Range("I1").Value = "FIRST AVE"
$red.Range("J1").Value = "ST SQUARE"
$red.Range("K1").Value = "STH HILLS"

# copy the bold/fill/center header formatting across the newly-used K column
$red.Range("J1").Copy()
$red.Range("K1").PasteSpecial(-4122)  # xlPasteFormats
$red.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3. "Red" sheet — driver + departure-time-string data
# ---------------------------------------------------------------------
$red.Range("B2").Value = "Employee A"
$red.Range("B3").Value = "Employee B"
$red.Range("B4").Value = "Employee C"
$red.Range("B5").Value = "Employee D"
$red.Range("B6").Value = "Employee E"

$red.Range("C2").Value = "06.00.00"
$red.Range("C3").Value = "06.15.00"
$red.Range("C4").Value = "06.30.00"
$red.Range("C5").Value = "06.45.00"
$red.Range("C6").Value = "07.00.00"

# old row 7 (train 105 / "Red") is gone - only a blank, General-formatted C7 remains
$red.Range("A7").Clear()
$red.Range("B7").Clear()
$red.Range("C7").NumberFormat = "General"

$red.Columns("B").ColumnWidth = 12.3

# ---------------------------------------------------------------------
# 4. "Green" sheet — header row (copy formatting from Red's header)
# ---------------------------------------------------------------------
$red.Range("A1:K1").Copy()
$green.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
$green.Application.CutCopyMode = $false

$green.Range("A1").Value = "Train ID"
$green.Range("B1").Value = "DRIVER"
$green.Range("C1").Value = "DEPARTURE"
$green.Range("D1").Value = "SHADYSIDE"
$green.Range("E1").Value = "HERRON"
$green.Range("F1").Value = "SWISSVILLE"
$green.Range("G1").Value = "PENN STATION"
$green.Range("H1").Value = "STEEL PLAZA"
$green.Range("I1").Value = "FIRST AVE"
$green.Range("J1").Value = "ST SQUARE"
$green.Range("K1").Value = "STH HILLS"

# ---------------------------------------------------------------------
# 5. "Green" sheet — driver + numeric departure times
# ---------------------------------------------------------------------
$green.Range("A2").Value = 100
$green.Range("A3").Value = 101
$green.Range("A4").Value = 102
$green.Range("A5").Value = 103
$green.Range("A6").Value = 104

$green.Range("B2").Value = "Employee A"
$green.Range("B3").Value = "Employee B"
$green.Range("B4").Value = "Employee C"
$green.Range("B5").Value = "Employee D"
$green.Range("B6").Value = "Employee E"

$green.Range("C2").Value = 0.25
$green.Range("C2").NumberFormat = "h:mm AM/PM"
$green.Range("C3").Value = 0.25486111111111109
$green.Range("C3").NumberFormat = "h:mm AM/PM"
$green.Range("C4").Value = 0.26041666666666669
$green.Range("C4").NumberFormat = "h:mm AM/PM"
$green.Range("C5").Value = 0.26527777777777778
$green.Range("C5").NumberFormat = "h:mm AM/PM"
$green.Range("C6").Value = 0.27083333333333331
$green.Range("C6").NumberFormat = "h:mm AM/PM"

# ---------------------------------------------------------------------
# 6. "Green" sheet — arrival-time schedule formulas (rows 7-11)
# ---------------------------------------------------------------------
$green.Range("D7").Formula = "=C2+TIME(0,3.7,0)"
$green.Range("E7").Formula = "=D7+TIME(0,2.3,0)"
$green.Range("F7").Formula = "=E7+TIME(0,1.5,0)"
$green.Range("G7").Formula = "=F7+TIME(0,1.8,0)"
$green.Range("H7").Formula = "=G7+TIME(0,2.1,0)"
$green.Range("I7").Formula = "=H7+TIME(0,2.1,0)"
$green.Range("J7").Formula = "=I7+TIME(0,1.7,0)"
$green.Range("K7").Formula = "=J7+TIME(0,2.3,0)"

$green.Range("D8").Formula = "=C3+TIME(0,3.7,0)"
$green.Range("E8").Formula = "=D8+TIME(0,2.3,0)"
$green.Range("F8").Formula = "=E8+TIME(0,1.5,0)"
$green.Range("G8").Formula = "=F8+TIME(0,1.8,0)"
$green.Range("H8").Formula = "=G8+TIME(0,2.1,0)"
$green.Range("I8").Formula = "=H8+TIME(0,2.1,0)"
$green.Range("J8").Formula = "=I8+TIME(0,1.7,0)"
$green.Range("K8").Formula = "=J8+TIME(0,2.3,0)"

$green.Range("D9").Formula = "=C4+TIME(0,3.7,0)"
$green.Range("E9").Formula = "=D9+TIME(0,2.3,0)"
$green.Range("F9").Formula = "=E9+TIME(0,1.5,0)"
$green.Range("G9").Formula = "=F9+TIME(0,1.8,0)"
$green.Range("H9").Formula = "=G9+TIME(0,2.1,0)"
$green.Range("I9").Formula = "=H9+TIME(0,2.1,0)"
$green.Range("J9").Formula = "=I9+TIME(0,1.7,0)"
$green.Range("K9").Formula = "=J9+TIME(0,2.3,0)"

$green.Range("D10").Formula = "=C5+TIME(0,3.7,0)"
$green.Range("E10").Formula = "=D10+TIME(0,2.3,0)"
$green.Range("F10").Formula = "=E10+TIME(0,1.5,0)"
$green.Range("G10").Formula = "=F10+TIME(0,1.8,0)"
$green.Range("H10").Formula = "=G10+TIME(0,2.1,0)"
$green.Range("I10").Formula = "=H10+TIME(0,2.1,0)"
$green.Range("J10").Formula = "=I10+TIME(0,1.7,0)"
$green.Range("K10").Formula = "=J10+TIME(0,2.3,0)"

$green.Range("D11").Formula = "=C6+TIME(0,3.7,0)"
$green.Range("E11").Formula = "=D11+TIME(0,2.3,0)"
$green.Range("F11").Formula = "=E11+TIME(0,1.5,0)"
$green.Range("G11").Formula = "=F11+TIME(0,1.8,0)"
$green.Range("H11").Formula = "=G11+TIME(0,2.1,0)"
$green.Range("I11").Formula = "=H11+TIME(0,2.1,0)"
$green.Range("J11").Formula = "=I11+TIME(0,1.7,0)"
$green.Range("K11").Formula = "=J11+TIME(0,2.3,0)"

$green.Range("D7:K11").NumberFormat = "h:mm AM/PM"

# ---------------------------------------------------------------------
# 7. Sheet views / selections
# ---------------------------------------------------------------------
$green.Range("F7").Select()

$red.Range("I2").Select()
$red.Activate()

